$wb = $excel.ActiveWorkbook

# --- Sheet "Month": just select B2 (cursor moved there) ---
$wsMonth = $wb.Worksheets.Item("Month")
$wsMonth.Range("B2").Select()

# --- Sheet "Year": just select B1 ---
$wsYear = $wb.Worksheets.Item("Year")
$wsYear.Range("B1").Select()

# --- Sheet "Day": scroll so row 10 is at top, select B26 ---
$wsDay = $wb.Worksheets.Item("Day")
$wsDay.Activate()
$excel.ActiveWindow.ScrollRow = 10
$wsDay.Range("B26").Select()

# --- Sheet "Hour": rename header "Sleepy" -> "Active", flip probabilities (1 - old), select A2 ---
$wsHour = $wb.Worksheets.Item("Hour")
$wsHour.Activate()

$wsHour.Range("B1").Value = "Active"

$wsHour.Range("B2").Value = 0.1
$wsHour.Range("B3").Value = 0.1
$wsHour.Range("B4").Value = 0.09999999999999998
$wsHour.Range("B5").Value = 0.09999999999999998
$wsHour.Range("B6").Value = 0.09999999999999998
$wsHour.Range("B7").Value = 0.09999999999999998
$wsHour.Range("B8").Value = 0.09999999999999998
$wsHour.Range("B9").Value = 0.30000000000000004
$wsHour.Range("B10").Value = 0.7
$wsHour.Range("B11").Value = 0.9
$wsHour.Range("B12").Value = 0.9
$wsHour.Range("B13").Value = 0.9
$wsHour.Range("B14").Value = 0.7
$wsHour.Range("B15").Value = 0.7
$wsHour.Range("B16").Value = 0.7
$wsHour.Range("B17").Value = 0.7
$wsHour.Range("B18").Value = 0.9
$wsHour.Range("B19").Value = 0.9
$wsHour.Range("B20").Value = 0.9
$wsHour.Range("B21").Value = 0.9
$wsHour.Range("B22").Value = 0.7
$wsHour.Range("B23").Value = 0.30000000000000004
$wsHour.Range("B24").Value = 0.09999999999999998

$wsHour.Range("A2").Select()
